# Generate Report for Handoff
# Replaces the source-file GUID token (fce084eb-...) with the newly
# generated one (555fed5e-...) across the "Overview", "zh-cn" and "de-de"
# sheets, and refreshes the handoff file names / handoff timestamps that
# depend on it.

$wb = $excel.ActiveWorkbook

$oldGuid = "fce084eb-eaca-4dcb-bbdd-d476fc9a2cf2"
$newGuid = "555fed5e-c795-47bb-ac58-0a5b4e106778"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldZhXlf = "$oldGuid.728838bf1020f2dd6693784a427fca9e1a6ca36c.zh-cn.xlf"
$newZhXlf = "$newGuid.775b2ed9aa23f58d66295b5b0fc81ae89f25cec6.zh-cn.xlf"

$oldDeXlf = "$oldGuid.728838bf1020f2dd6693784a427fca9e1a6ca36c.de-de.xlf"
$newDeXlf = "$newGuid.775b2ed9aa23f58d66295b5b0fc81ae89f25cec6.de-de.xlf"

$oldZhTime = "2016-03-11 01:27:07"
$newZhTime = "2016-03-11 01:27:48"

$oldDeTime = "2016-03-11 01:27:14"
$newDeTime = "2016-03-11 01:27:54"

# Original external hyperlink targets (unchanged by this edit - only the
# visible "display" text changes in the diff) so the relationship keeps
# pointing at the right place after we refresh the display text.
$mdTargetOverview = "https://github.com/OpenLocalizationTest/oltest/blob/ca78077d86b08d918e9cd9f9c5bc70f317a45ff3/e2e/$oldMd"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ab80e6c14f26fb9c58b11a2a4bbe95bd5447c2a2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a44a44c8b97fc77317a3405107d80ba151d51165/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

function Update-HyperlinkDisplay($ws, $cellRef, $newDisplay, $targetUrl) {
    # Find the hyperlink currently attached to this cell.
    $range = $ws.Range($cellRef)
    $count = $range.Hyperlinks.Count
    if ($count -ge 1) {
        $hl = $range.Hyperlinks.Item(1)
        # Updates display text (appends a corrected overlay entry under the
        # hood); immediately restore the Address so the relationship/target
        # survives the refresh instead of being dropped.
        $hl.TextToDisplay = $newDisplay
        $refreshed = $ws.Hyperlinks.Item($ws.Hyperlinks.Count)
        $refreshed.Address = $targetUrl
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview": A2 holds the handoff markdown file name/hyperlink.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
Update-HyperlinkDisplay $wsOverview "A2" $newMd $mdTargetOverview

# ---------------------------------------------------------------------
# Sheet "zh-cn": A2 markdown file, C2 handoff xlf file, D2 handoff time.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
Update-HyperlinkDisplay $wsZh "A2" $newMd $mdTargetOverview

$wsZh.Range("C2").Value = $newZhXlf
Update-HyperlinkDisplay $wsZh "C2" $newZhXlf $zhXlfTarget

$wsZh.Range("D2").Value = $newZhTime

# ---------------------------------------------------------------------
# Sheet "de-de": A2 markdown file, C2 handoff xlf file, D2 handoff time.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
Update-HyperlinkDisplay $wsDe "A2" $newMd $mdTargetOverview

$wsDe.Range("C2").Value = $newDeXlf
Update-HyperlinkDisplay $wsDe "C2" $newDeXlf $deXlfTarget

$wsDe.Range("D2").Value = $newDeTime
